# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. This updates the DAMSLTag (column I) and
# DialogAct (column J) values for a number of rows on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Row = 3;   Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 19;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 25;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 26;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 32;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 55;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 62;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 79;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 80;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 82;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 99;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 100; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 104; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 105; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 112; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 126; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 129; Tag = "%";  Act = "Uninterpretable" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.Tag   # Column I - DAMSLTag
    $ws.Cells.Item($u.Row, 10).Value = $u.Act  # Column J - DialogAct
}
